$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "[name=`"Talulah?`"] Or did you think that I would be suppressed by the dragon-slaying sword of Yen and the filthy Arts of the Sarkaz? Free for your Arts to do with this body as you wish, to regain some semblance of 'purity'?`n"
$ws.Range("C46").Value = "[name=`"Talulah?`"] 'I journey across the plains of salt, the bleeding mountain valleys, and not once was I harmed.' `n"
$ws.Range("C47").Value = "[name=`"Talulah?`"] 'For his sword, his territory, his subjects, they grant me protection.' `n"
$ws.Range("C114").Value = "[name=`"'Deathless Black Snake'`"] Amiya, oh Amiya. Do you know why you are so shortsighted, why you are so weak?`n"
$ws.Range("C115").Value = "[name=`"'Deathless Black Snake'`"] Because you only stand in your own two shoes. Do you know how narrow your vision is? How much misfortune have you ever seen?`n"
$ws.Range("C116").Value = "[name=`"'Deathless Black Snake'`"] You’ve only seen some unremarkable scenes through the eyes of the Lord of Fiends, yet you presume to act as if you have personally experienced it...`n"
$ws.Range("C117").Value = "[name=`"'Deathless Black Snake'`"] I’ve experienced far more tragedy than you ever have.`n"
$ws.Range("C118").Value = "[name=`"'Deathless Black Snake'`"] You are constrained by your identity as an orphaned Infected, whereas I am the observer of the entirety of Ursus.`n"
$ws.Range("C119").Value = "[name=`"'Deathless Black Snake'`"] I live to save the vast land and crying people you cannot see. What about you?`n"
$ws.Range("C120").Value = "[name=`"'Deathless Black Snake'`"] You are but fools, blind and deaf, despairing over the littlest of things, proclaiming yourself saviors. `n"
$ws.Range("C121").Value = "[name=`"'Deathless Black Snake'`"] Could you ever imagine what has taken place across this world? You cannot.`n"
$ws.Range("C122").Value = "[name=`"'Deathless Black Snake'`"] The language, history and memories you see are all just short excerpts of self-glorification. You cannot see the scars carved onto the living by fate itself.`n"
$ws.Range("C123").Value = "[name=`"'Deathless Black Snake'`"] Scars can kill. The scars of millions of my Ursus kin are etched onto my body.`n"
$ws.Range("C124").Value = "[name=`"'Deathless Black Snake'`"] Ch'en Hui-chieh. You are my sister in blood and flesh. You stand here against me and view Ursus as your enemy, all because you think it is invading your homeland.`n"
$ws.Range("C125").Value = "[name=`"'Deathless Black Snake'`"] If you cannot even defeat your own sister, do you believe for a second that you will prevail against the entirety of Ursus?`n"
$ws.Range("C126").Value = "[name=`"'Deathless Black Snake'`"] What’s more, whether or not Ursus will be dragged into this war, and whatever the result of the war may be, matters not to me.`n"
$ws.Range("C127").Value = "[name=`"'Deathless Black Snake'`"] All I want is for this land to free itself from the impasse and dissonance it finds itself in. After this war ends, there will inevitably be a winner and a loser, whether it is the military government or a new emperor.`n"
$ws.Range("C128").Value = "[name=`"'Deathless Black Snake'`"] The previous Ursus emperor was like the incarnation of the soul of Ursus.`n"
$ws.Range("C129").Value = "[name=`"'Deathless Black Snake'`"] He did not lose himself in pleasure, and maintained a healthy distance from power despite his lofty ambitions. He could not be called the embodiment of lust, and was as far from a tyrant as could be.`n"
$ws.Range("C130").Value = "[name=`"'Deathless Black Snake'`"] However, he brought about wars, for Ursus needs wars. Ursus needs wars for its people to grow.`n"
$ws.Range("C131").Value = "[name=`"'Deathless Black Snake'`"] He was exceptional. He was worthy of being called an emperor of Ursus.`n"
$ws.Range("C132").Value = "[name=`"'Deathless Black Snake'`"] As for you, Infected.`n"
$ws.Range("C133").Value = "[name=`"'Deathless Black Snake'`"] 'Why don’t you bow down in front of them and wag your tails desperately like the beasts you are? You may receive their forgiveness that way.' `n"
$ws.Range("C134").Value = "[name=`"'Deathless Black Snake'`"] You wish to be accepted by everyone?`n"
$ws.Range("C135").Value = "[name=`"'Deathless Black Snake'`"] Of course you can... Just accept the fodder they give you, abandon the language of people and become slaves. Hand over a quarter of your own for them to squander and slaughter.`n"
$ws.Range("C136").Value = "[name=`"'Deathless Black Snake'`"] They will see you as harmless livestock, a pet that looks disgustingly like a person yet is interesting enough to keep around. An accessory that they may grant meager protection to console their miserable pride.`n"
$ws.Range("C137").Value = "[name=`"'Deathless Black Snake'`"] You will struggle and die, yet there will always be countless new lives ready to take your place. Your kin will be reduced to such, then and forever.`n"
$ws.Range("C138").Value = "[name=`"'Deathless Black Snake'`"] ——Because those fools can only feel forgiveness and satisfaction through such a method. The great systems of power are always established for the sake of you and them... the slave and the master.`n"
$ws.Range("C139").Value = "[name=`"'Deathless Black Snake'`"] No one else will treat you as fairly as Ursus and I. We know how despicable the Ursus may be, yet we have no choice but to let you live in their shadow.`n"
$ws.Range("C140").Value = "[name=`"'Deathless Black Snake'`"] If you think that is evil, then so be it. Because this evil perpetuated a mentality of pathetic salvation that spanned nations, that caused life to be so.`n"
$ws.Range("C141").Value = "[name=`"'Deathless Black Snake'`"] Until today.`n"
$ws.Range("C142").Value = "[name=`"'Deathless Black Snake'`"] I will change everything. I gave them direction.`n"
$ws.Range("C143").Value = "[name=`"'Deathless Black Snake'`"] The exceptional once-ruler of this nation brought greatness to them, yet even he could not escape the grasp of time. In the end, he could not bring salvation to the Infected.`n"
$ws.Range("C144").Value = "[name=`"'Deathless Black Snake'`"] In war, their irrational parts are sated, their non-intelligent parts are perfused, their mediocre parts become extraordinary, their indifferent thoughts become fanatic.`n"
$ws.Range("C145").Value = "[name=`"'Deathless Black Snake'`"] Those people gain a sort of nobility.`n"
$ws.Range("C146").Value = "[name=`"'Deathless Black Snake'`"] To tame and use them is an act of foolishness rather than wisdom. We do not need their allegiance, nor do we want to rule over them.`n"
$ws.Range("C147").Value = "[name=`"'Deathless Black Snake'`"] I will let them restore their lost dignity.`n"
$ws.Range("C148").Value = "[name=`"'Deathless Black Snake'`"] In the future of Ursus, the same holds true for the Infected.`n"
$ws.Range("C149").Value = "[name=`"'Deathless Black Snake'`"] I will bring about exaltation and equality, for the Infected and non-Infected, for Ursus and other nations.`n"
$ws.Range("C150").Value = "[name=`"'Deathless Black Snake'`"] Inevitably, a winner will emerge between the borderland armies and the State Council. I fan the flames of this war, so that they may decide immediately.`n"
$ws.Range("C151").Value = "[name=`"'Deathless Black Snake'`"] As time passes, life is eroded and passion is ousted. When Ursus becomes incapable of setting off the slightest spark, its lands and people will rot and become lost in the annals of history.`n"
$ws.Range("C152").Value = "[name=`"'Deathless Black Snake'`"] I will never let that happen.`n"
$ws.Range("C153").Value = "[name=`"'Deathless Black Snake'`"] I had intentionally let Reunion fall here, to be regathered in my hand and crumble once again, just so that there is a starting point to all this.`n"
$ws.Range("C154").Value = "[name=`"'Deathless Black Snake'`"] Yet, you deny all of this. You deny everything that I will do for Ursus.`n"
$ws.Range("C155").Value = "[name=`"'Deathless Black Snake'`"] How can your shortsightedness see the events a century henceforth?`n"
$ws.Range("C156").Value = "[name=`"'Deathless Black Snake'`"] How can you even begin to imagine the Infected and non-Infected united in their struggle? Why would you disrupt the peace and glory that will follow a thousand years, simply for the sake of a little suffering today?`n"
$ws.Range("C157").Value = "[name=`"'Deathless Black Snake'`"] In the end, you are too immature.`n"
$ws.Range("C158").Value = "[name=`"'Deathless Black Snake'`"] You have never experienced despair at the brink of destruction. The Infected may struggle, but they still live on. In contrast, Ursus has nearly fallen to ruin in a single night on multiple occasions.`n"
$ws.Range("C159").Value = "[name=`"'Deathless Black Snake'`"] Have you ever listened to the cavalry of Khagan galloping across the land, the sound of weapons clashing so thunderously they can level mountains and shift rivers?`n"
$ws.Range("C160").Value = "[name=`"'Deathless Black Snake'`"] Have you ever heard the shrill cries of the Kheshigs of Nightzmora as they run their scimitars across their skulls?`n"
$ws.Range("C161").Value = "[name=`"'Deathless Black Snake'`"] Have you ever experienced the ear-splitting sound of thousands of Warriors of Rus firing simultaneously, reducing the high and mighty Gauls to little more than ashes amongst the blood and mud?`n"
$ws.Range("C162").Value = "[name=`"'Deathless Black Snake'`"] You are clueless, but I don’t blame you for that.`n"
$ws.Range("C163").Value = "[name=`"'Deathless Black Snake'`"] If you really wish to put your foolish thoughts into action and prevent the revival of Ursus, the reunification of the Ursus under the same banner of faith, the revitalization of the Ursus land...`n"
$ws.Range("C164").Value = "[name=`"'Deathless Black Snake'`"] Then, try and stop me. Try and stop this incarnation of Ursus.`n"
$ws.Range("C165").Value = "[name=`"'Deathless Black Snake'`"] Come.`n"
$ws.Range("D3").Value = "[name=`"탈룰라?`"] 하지만 너와 네 곁에 있는 어린 용의 헛된 자신감이 꽤 불쾌하구나.`n"
$ws.Range("D36").Value = "아미야가 첸의 머릿속의 이미지를 자신의 머릿속으로 보내 달라고 손짓했다. 첸은 깨달았다, 이것은 분명 아츠였다.`n"
